$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 312.53845
$ws.Range("I5").Value = 96.2
$ws.Range("J5").Value = 447.75
$ws.Range("K5").Value = 96.2
$ws.Range("L5").Value = 447.75
$ws.Range("M5").Value = 18.8
$ws.Range("N5").Value = -677.75
$ws.Range("H33").Value = 289.25
$ws.Range("I33").Value = 254.73685
$ws.Range("J33").Value = 420.4
$ws.Range("K33").Value = 254.73685
$ws.Range("L33").Value = 420.4
$ws.Range("M33").Value = -25.73685
$ws.Range("N33").Value = -878.4
$ws.Range("H40").Value = 1355.9706
$ws.Range("I40").Value = 1351.625
$ws.Range("J40").Value = 1366.4
$ws.Range("K40").Value = 1351.625
$ws.Range("L40").Value = 1366.4
$ws.Range("M40").Value = -1176.625
$ws.Range("N40").Value = -1716.4
$ws.Range("H64").Value = 3714.8108
$ws.Range("I64").Value = 3528
$ws.Range("J64").Value = 3960
$ws.Range("K64").Value = 3528
$ws.Range("L64").Value = 3960
$ws.Range("M64").Value = -3280
$ws.Range("N64").Value = -4456
$ws.Range("H67").Value = 3714.8108
$ws.Range("I67").Value = 3528
$ws.Range("J67").Value = 3960
$ws.Range("K67").Value = 3528
$ws.Range("L67").Value = 3960
$ws.Range("M67").Value = -2670
$ws.Range("N67").Value = -5676
$ws.Range("H74").Value = 4327.143
$ws.Range("J74").Value = 4418.2
$ws.Range("L74").Value = 4418.2
$ws.Range("N74").Value = -6290.2
$ws.Range("H77").Value = 4327.143
$ws.Range("J77").Value = 4418.2
$ws.Range("L77").Value = 22091
$ws.Range("N77").Value = -31451
$ws.Range("H86").Value = 2889.2222
$ws.Range("I86").Value = 1901
$ws.Range("J86").Value = 3383.3333
$ws.Range("K86").Value = 1901
$ws.Range("L86").Value = 3383.3333
$ws.Range("M86").Value = -778
$ws.Range("N86").Value = -5629.3333
$ws.Range("H89").Value = 2889.2222
$ws.Range("I89").Value = 1901
$ws.Range("J89").Value = 3383.3333
$ws.Range("K89").Value = 9505
$ws.Range("L89").Value = 16916.6665
$ws.Range("M89").Value = -3889
$ws.Range("N89").Value = -28148.6665
$ws.Range("H106").Value = 2777.7778
$ws.Range("I106").Value = 2000
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 2000
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -1369
$ws.Range("N106").Value = -4262

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1278
$ws.Range("I2").Value = 1329
$ws.Range("J2").Value = 1099.5
$ws.Range("K2").Value = 1329
$ws.Range("L2").Value = 1099.5
$ws.Range("M2").Value = -1216
$ws.Range("N2").Value = -1325.5
$ws.Range("H116").Value = 1278
$ws.Range("I116").Value = 1329
$ws.Range("J116").Value = 1099.5
$ws.Range("K116").Value = 1329
$ws.Range("L116").Value = 1099.5
$ws.Range("M116").Value = 965
$ws.Range("N116").Value = -5687.5
$ws.Range("H132").Value = 3501.738
$ws.Range("I132").Value = 2542.12
$ws.Range("J132").Value = 4912.9414
$ws.Range("K132").Value = 7626.36
$ws.Range("L132").Value = 14738.8242
$ws.Range("M132").Value = -5096.36
$ws.Range("N132").Value = -19798.8242

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1278
$ws.Range("I3").Value = 1329
$ws.Range("J3").Value = 1099.5
$ws.Range("K3").Value = 1329
$ws.Range("L3").Value = 1099.5
$ws.Range("M3").Value = -1215
$ws.Range("N3").Value = -1327.5
$ws.Range("H134").Value = 2282.7827
$ws.Range("I134").Value = 1536
$ws.Range("K134").Value = 4608
$ws.Range("M134").Value = -2073

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3753.7273
$ws.Range("I122").Value = 3657.2856
$ws.Range("J122").Value = 3922.5
$ws.Range("K122").Value = 10971.8568
$ws.Range("L122").Value = 11767.5
$ws.Range("M122").Value = -8521.856800000001
$ws.Range("N122").Value = -16667.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4507958
$ws.Range("I5").Value = 421.73914
$ws.Range("J5").Value = 11913197
$ws.Range("K5").Value = 1265.21742
$ws.Range("L5").Value = 35739591
$ws.Range("M5").Value = -1153.21742
$ws.Range("N5").Value = -35739815
$ws.Range("H22").Value = 1517.0834
$ws.Range("J22").Value = 1817.4445
$ws.Range("L22").Value = 5452.333500000001
$ws.Range("N22").Value = -5790.333500000001
$ws.Range("H27").Value = 1517.0834
$ws.Range("J27").Value = 1817.4445
$ws.Range("L27").Value = 5452.333500000001
$ws.Range("N27").Value = -5656.333500000001
$ws.Range("H41").Value = 527.5
$ws.Range("I41").Value = 167
$ws.Range("J41").Value = 888
$ws.Range("K41").Value = 501
$ws.Range("L41").Value = 2664
$ws.Range("M41").Value = -163
$ws.Range("N41").Value = -3340
$ws.Range("H97").Value = 12045
$ws.Range("I97").Value = 1064.2858
$ws.Range("J97").Value = 37666.668
$ws.Range("K97").Value = 3192.8574
$ws.Range("L97").Value = 113000.004
$ws.Range("M97").Value = -2696.8574
$ws.Range("N97").Value = -113992.004
$ws.Range("H135").Value = 4507958
$ws.Range("I135").Value = 421.73914
$ws.Range("J135").Value = 11913197
$ws.Range("K135").Value = 3795.65226
$ws.Range("L135").Value = 107218773
$ws.Range("M135").Value = -1260.65226
$ws.Range("N135").Value = -107223843

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 18333.334
$ws.Range("I80").Value = 25500
$ws.Range("K80").Value = 25500
$ws.Range("M80").Value = -24502
$ws.Range("H83").Value = 18333.334
$ws.Range("I83").Value = 25500
$ws.Range("K83").Value = 127500
$ws.Range("M83").Value = -122508

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4644.2383
$ws.Range("I7").Value = 3582.6155
$ws.Range("J7").Value = 6369.375
$ws.Range("K7").Value = 3582.6155
$ws.Range("L7").Value = 6369.375
$ws.Range("M7").Value = -3470.6155
$ws.Range("N7").Value = -6593.375
$ws.Range("H126").Value = 4644.2383
$ws.Range("I126").Value = 3582.6155
$ws.Range("J126").Value = 6369.375
$ws.Range("K126").Value = 10747.8465
$ws.Range("L126").Value = 19108.125
$ws.Range("M126").Value = -8277.8465
$ws.Range("N126").Value = -24048.125
$ws.Range("H131").Value = 49714
$ws.Range("J131").Value = 49714
$ws.Range("L131").Value = 49714
$ws.Range("N131").Value = -59794
$ws.Range("H132").Value = 3940.2942
$ws.Range("I132").Value = 3325.8462
$ws.Range("J132").Value = 5937.25
$ws.Range("K132").Value = 9977.5386
$ws.Range("L132").Value = 17811.75
$ws.Range("M132").Value = -7447.5386
$ws.Range("N132").Value = -22871.75
